# "Generate Report for Archive"
# - Update status text from "Ready for handoff" to "In Translation" on all
#   three sheets (Overview, zh-cn, de-de).
# - Narrow the now-shorter "Status" columns to fit the new text
#   (Overview!E:F and the Status column (C) on the zh-cn / de-de sheets).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Replace the "Ready for handoff" status text with "In Translation"
# everywhere it appears.
$wsOverview.Range("E2:F4").Value = "In Translation"
$wsZhCn.Range("C2:C4").Value = "In Translation"
$wsDeDe.Range("C2:C4").Value = "In Translation"

# Re-fit the Status columns now that the text is shorter.
$wsOverview.Range("E1:F4").ColumnWidth = 12.5
$wsZhCn.Range("C1:C4").ColumnWidth = 12.5
$wsDeDe.Range("C1:C4").ColumnWidth = 12.5
